$wb = $excel.ActiveWorkbook

# Map of cell address -> new value for the "想去人数" (attendance) column.
$updates = @{
    "F3"  = 2982
    "F7"  = 1630
    "F9"  = 79
    "F11" = 1337
    "F13" = 468
    "F19" = 102
    "F20" = 3077
    "F21" = 374
    "F22" = 95
    "F24" = 87
}

# Both "展览" (sheet1) and "全部类型" (sheet4) contain the same rows of data
# and both need to be updated with the refreshed counts.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
